$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "67.905.47"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "3.781.82"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.86"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.08"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("D7").Value = "3.778.00"
$ws.Range("E7").Value = "  -1.25%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  -2.37%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.446"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.81"
$ws.Range("E12").Value = "  +8.11%  "
$ws.Range("E13").Value = "  -2.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.99"
$ws.Range("E14").Value = "  -2.65%  "
$ws.Range("D15").Value = "4.418.53"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Value = "3.786.11"
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("D17").Value = "67.878.86"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("E19").Value = "  +1.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.99"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "459.55"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.44"
$ws.Range("E22").Value = "  -4.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.689"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("E25").Value = "  -4.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.85"
$ws.Range("E26").Value = "  -2.05%  "
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.89"
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("D30").Value = "3.938.35"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.21"
$ws.Range("E31").Value = "  -2.61%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.59"
$ws.Range("E32").Value = "  -7.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.18"
$ws.Range("E33").Value = "  -2.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.95"
$ws.Range("E34").Value = "  -2.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.90"
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0990"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("E38").Value = "  +5.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.80"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("E40").Value = "  -4.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.976"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.65"
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.09"
$ws.Range("E45").Value = "  -2.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "152.27"
$ws.Range("E46").Value = "  +2.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.293"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("E48").Value = "  -2.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.27"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.82"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.40"
$ws.Range("E51").Value = "  -9.93%  "
